$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.627.24"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "3.368.04"
$ws.Range("E3").Value = "  +0.52%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "562.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "176.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.34%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.73%  "

$ws.Range("D8").Value = "3.358.95"
$ws.Range("E8").Value = "  +0.47%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.635"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.79%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +10.73%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "55.81"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.77%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000277"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.70%  "

$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("D15").Value = "3.900.84"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("E16").Value = "  +4.23%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.366.81"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  +3.05%  "

$ws.Range("D20").Value = "64.507.05"
$ws.Range("E20").Value = "  +2.15%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.79%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "462.53"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +15.43%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +13.30%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.92%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "86.45"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.97%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "13.59"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.84%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.86"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "

$ws.Range("E28").Value = "  +5.35%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.86"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("E30").Value = "  +4.30%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.68"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.99%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "11.53"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.64%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "582.38"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("E34").Value = "  +3.45%  "

$ws.Range("E35").Value = "  +2.64%  "

$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("E37").Value = "  -5.82%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "36.04"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("D39").Value = "0.0₃0760"
$ws.Range("E39").Value = "  +5.15%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.46"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.372"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "

$ws.Range("D42").Value = "3.096.76"
$ws.Range("E42").Value = "  -1.56%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("E44").Value = "  +1.12%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("E46").Value = "  +3.94%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.133"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.34%  "

$ws.Range("E49").Value = "  +0.46%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.44"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.84%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "137.23"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.90%  "
